$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c_E2 = @'

                    <filter>
                        <interfaces xmlns="http://openconfig.net/yang/interfaces">
                        <interface>
                        <name>lag-1</name>
                        </interface>
                        </interfaces>
                    </filter>
                     
-------------------

                    <filter>
                        <interfaces xmlns="http://openconfig.net/yang/interfaces">
                        <interface>
                        <name>1/1/1</name>
                        </interface>
                        </interfaces>
                    </filter>
                    
'@
$ws.Range("E2").Value = $c_E2

$c_F2 = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:746e106d-83da-4bfe-8f4b-9bff3c5d2e3d" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
    </data>
</rpc-reply> 
-------------------
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:dde397d5-a583-4a97-b71f-7b1cd1f29e58" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <auto-negotiate>false</auto-negotiate>
                        <duplex-mode>FULL</duplex-mode>
                        <port-speed>SPEED_100MB</port-speed>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>
'@
$ws.Range("F2").Value = $c_F2

$c_G2 = @'
  <edit-config>
      <target>
        <candidate/>
      </target>
      <config>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
          <interface>
            <name>lag-1</name>
            <config>
              <name>lag-1</name>
              <type>ieee8023adLag</type>
              <enabled>true</enabled>
            </config>
          </interface>
        </interfaces>
      </config>
    </edit-config> 
-------------------
<edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <interfaces xmlns="http://openconfig.net/yang/interfaces">
        <interface>
          <name>1/1/1</name>
          <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
            <config>
              <aggregate-id xmlns="http://openconfig.net/yang/interfaces/aggregate">lag-1</aggregate-id>
            </config>
          </ethernet>
        </interface>
      </interfaces>
    </config>
  </edit-config>
'@
$ws.Range("G2").Value = $c_G2

$c_H2 = @'
- Response of edit-config: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:8331c687-afda-4095-9577-1643acc83fa0" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply> 

 - Response of commit: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:f657f7cb-0ab8-49e4-a995-76e30e2a2bec" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply> 
-------------------
- Response of edit-config: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:130c8d38-d7bb-4460-8a01-3cf030f783e7" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply> 

 - Response of commit: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:e1d22199-55bf-4b43-8311-80f09b72deb3" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply>
'@
$ws.Range("H2").Value = $c_H2

$c_I2 = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:c2d6484d-7f1b-4e11-a381-56e67adbf4bd" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>lag-1</name>
                <config>
                    <name>lag-1</name>
                    <type>ieee8023adLag</type>
                    <enabled>true</enabled>
                </config>
            </interface>
        </interfaces>
    </data>
</rpc-reply> 
-------------------
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:2d61ec01-e33f-4989-a66f-ca004e7426c6" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <auto-negotiate>false</auto-negotiate>
                        <duplex-mode>FULL</duplex-mode>
                        <port-speed>SPEED_100MB</port-speed>
                        <aggregate-id>lag-1</aggregate-id>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>
'@
$ws.Range("I2").Value = $c_I2

